$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 270, pushing the existing data (old rows 270-294)
# down to rows 272-296.
$ws.Rows.Item(270).Resize(2).Insert()

# ---- New row 270: Lapins / Primera, week of 2022-01-24 ----
$ws.Cells.Item(270, 1).Value2  = 8
$ws.Cells.Item(270, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(270, 3).Value2  = "Coquimbo"
$ws.Cells.Item(270, 4).Value2  = 44585
$ws.Cells.Item(270, 5).Value2  = 4
$ws.Cells.Item(270, 6).Value2  = "Fruta"
$ws.Cells.Item(270, 7).Value2  = 100103
$ws.Cells.Item(270, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(270, 9).Value2  = 100103001
$ws.Cells.Item(270, 10).Value2 = "Cereza"
$ws.Cells.Item(270, 11).Value2 = "Lapins"
$ws.Cells.Item(270, 12).Value2 = "Primera"
$ws.Cells.Item(270, 13).Value2 = 400
$ws.Cells.Item(270, 14).Value2 = 9500
$ws.Cells.Item(270, 15).Value2 = 10000
$ws.Cells.Item(270, 16).Value2 = 9750
$ws.Cells.Item(270, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(270, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(270, 19).Value2 = 975
$ws.Cells.Item(270, 20).Value2 = 10

# ---- New row 271: Lapins / Segunda, week of 2022-01-24 ----
$ws.Cells.Item(271, 1).Value2  = 8
$ws.Cells.Item(271, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(271, 3).Value2  = "Coquimbo"
$ws.Cells.Item(271, 4).Value2  = 44585
$ws.Cells.Item(271, 5).Value2  = 4
$ws.Cells.Item(271, 6).Value2  = "Fruta"
$ws.Cells.Item(271, 7).Value2  = 100103
$ws.Cells.Item(271, 8).Value2  = "Frutos de hueso (carozo)"
$ws.Cells.Item(271, 9).Value2  = 100103001
$ws.Cells.Item(271, 10).Value2 = "Cereza"
$ws.Cells.Item(271, 11).Value2 = "Lapins"
$ws.Cells.Item(271, 12).Value2 = "Segunda"
$ws.Cells.Item(271, 13).Value2 = 500
$ws.Cells.Item(271, 14).Value2 = 7500
$ws.Cells.Item(271, 15).Value2 = 8000
$ws.Cells.Item(271, 16).Value2 = 7750
$ws.Cells.Item(271, 17).Value2 = "`$/bandeja 10 kilos"
$ws.Cells.Item(271, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(271, 19).Value2 = 775
$ws.Cells.Item(271, 20).Value2 = 10

$ws.Range("A1").Select()
